# Updates cryptos list values (Price / Volume(1h)) per "Updated cryptos list" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '74.598.59'
$ws.Range("E2").Value = '  +6.39%  '
# Row 3
$ws.Range("D3").Value = '2.658.61'
$ws.Range("E3").Value = '  +8.13%  '
# Row 4
$ws.Range("E4").Value = '  -0.05%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '186.23'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +10.94%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '583.41'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.48%  '
# Row 7
$ws.Range("E7").Value = '  -0.14%  '
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.533'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.71%  '
# Row 9
$ws.Range("E9").Value = '  +8.29%  '
# Row 10
$ws.Range("D10").Value = '2.658.43'
$ws.Range("E10").Value = '  +8.22%  '
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.164'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.26%  '
# Row 12
$ws.Range("E12").Value = '  +5.27%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.74'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.53%  '
# Row 14
$ws.Range("D14").Value = '3.144.32'
$ws.Range("E14").Value = '  +7.96%  '
# Row 15
$ws.Range("D15").Value = '74.508.39'
$ws.Range("E15").Value = '  +6.44%  '
# Row 16
$ws.Range("E16").Value = '  +1.33%  '
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.35'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +8.64%  '
# Row 18
$ws.Range("D18").Value = '2.658.18'
$ws.Range("E18").Value = '  +8.01%  '
# Row 19
$ws.Range("E19").Value = '  +29.09%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.86'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +8.96%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '371.11'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +8.03%  '
# Row 22
$ws.Range("E22").Value = '  +10.08%  '
# Row 23
$ws.Range("E23").Value = '  +4.08%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.24'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.16%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.11%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '69.22'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.79%  '
# Row 27
$ws.Range("E27").Value = '  +5.62%  '
# Row 28
$ws.Range("E28").Value = '  +8.35%  '
# Row 29
$ws.Range("D29").Value = '2.794.95'
$ws.Range("E29").Value = '  +7.92%  '
# Row 30
$ws.Range("E30").Value = '  +1.22%  '
# Row 31
$ws.Range("E31").Value = '  +7.95%  '
# Row 32
$ws.Range("E32").Value = '  +12.41%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '517.67'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +12.95%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.61'
$ws.Range("D34").Style = "Normal"
# Row 35
$ws.Range("E35").Value = '  +6.35%  '
# Row 36
$ws.Range("E36").Value = '  -0.15%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '162.66'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.24%  '
# Row 38
$ws.Range("E38").Value = '  +5.01%  '
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.20'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.06%  '
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '19.36'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.37%  '
# Row 41
$ws.Range("E41").Value = '  -0.03%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '169.25'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +25.81%  '
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.93'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +11.07%  '
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.327'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +6.93%  '
# Row 45
$ws.Range("E45").Value = '  +7.29%  '
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.18'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +7.11%  '
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '39.03'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.80%  '
# Row 48
$ws.Range("E48").Value = '  +8.53%  '
# Row 49
$ws.Range("E49").Value = '  +15.50%  '
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.62'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.98%  '
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '21.15'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +20.83%  '
